# Adding four new variables to spreadsheet
# Inserts 4 new rows at the top of the "New Variables" block (rows 176-179),
# describing tot_tran_in / tot_tran_out / tot_amount_in / tot_amount_out,
# formatted like the existing highlighted header rows (bold mono font,
# thin border, light gold fill) and pushes the remainder of the sheet down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room: insert 4 blank rows before the existing row 176.
$ws.Rows("176:179").Insert()

# 2. Populate the new rows - column A (variable name) and column B
#    (description). Write order matters for shared-string interning order,
#    so column A is written out of row-order (176,179,178,177) followed by
#    column B in row-order, matching how the strings were originally typed.
$ws.Range("A176").Value2 = "tot_tran_in                                    numeric"
$ws.Range("A179").Value2 = "tot_amount_out                                 numeric"
$ws.Range("A178").Value2 = "tot_amount_in                                  numeric"
$ws.Range("A177").Value2 = "tot_tran_out                                   numeric"

$ws.Range("B176").Value2 = "total transactions within offer periods"
$ws.Range("B177").Value2 = "total transactions not within offer periods"
$ws.Range("B178").Value2 = "total amount within offer periods"
$ws.Range("B179").Value2 = "total amount not within offer periods"

# 3. Format the new block (A:C, 176:179) to match the other highlighted
#    "variable name / description" rows elsewhere in the sheet: thin black
#    border around every cell, and a light gold (Accent4, 80% lighter) fill.
$newRows = $ws.Range("A176:C179")
$newRows.Borders.Color = 0
$newRows.Borders.Weight = 2
$newRows.Borders.LineStyle = 1
$newRows.Interior.Color = 13431551

# 4. Restore the user's selection near where they left off editing.
$ws.Range("B180").Select()
